$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Mon_Dec__4_07_13_32_2023"
$ws.Range("C1").Value = 30

# Row 2
$ws.Range("A2").Value = "Mon_Dec__4_07_17_55_2023"
$ws.Range("B2").Font.Bold = $false
$ws.Range("C2").Value = 30
